$d = $word.ActiveDocument

$replacements = @(
    @("2024-05-13 Monday", "2024-05-14 Tuesday"),
    @("491×6=2946", "121×2=242"),
    @("134×4=536", "829×6=4974"),
    @("740×8=5920", "703×5=3515"),
    @("686×8=5488", "985×8=7880"),
    @("952×9=8568", "952×2=1904"),
    @("550×2=1100", "392×9=3528"),
    @("305×3=915", "416×2=832"),
    @("273×4=1092", "438×9=3942"),
    @("617×8=4936", "708×8=5664"),
    @("750×6=4500", "341×3=1023"),
    @("646×8=5168", "323×9=2907"),
    @("580×3=1740", "511×4=2044"),
    @("808×5=4040", "904×7=6328"),
    @("840×8=6720", "998×5=4990"),
    @("105×4=420", "947×5=4735"),
    @("597×2=1194", "865×9=7785"),
    @("831×6=4986", "742×9=6678"),
    @("635×7=4445", "858×5=4290"),
    @("686×9=6174", "639×4=2556"),
    @("707×5=3535", "877×9=7893"),
    @("522×7=3654", "201×2=402"),
    @("838×3=2514", "366×3=1098"),
    @("561×7=3927", "491×7=3437"),
    @("679×4=2716", "883×4=3532"),
    @("678×6=4068", "605×7=4235")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
